$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.707.48'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.948.77'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.59%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.20'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.10'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +10.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.680'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.783'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.25'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000326'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.39%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.51'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.584.01'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.964.64'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.02'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.24%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.23'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.99'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.691.51'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.130'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '443.78'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.84'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '95.83'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.34'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -5.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.10'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.25'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.89'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('B29').Value = 'Filecoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.35'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.66'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.81'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '50.11'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.128'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0995'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +13.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '68.99'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '635.89'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -7.73%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.426'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.77%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.43'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.25%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.16%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.145'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.05%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0477'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('B44').Value = 'THORChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.52'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -7.44%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.13'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +41.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.147'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.79%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.39'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.61'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.87'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -14.35%  '
$ws.Range('B50').Value = 'FLOKI'
$ws.Range('C50').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000285'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.832.19'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.28%  '
